$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-28 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-01 Saturday", 2) | Out-Null
$d.Content.Find.Execute("67×37=2479", $true, $false, $false, $false, $false, $true, 1, $false, "24×83=1992", 2) | Out-Null
$d.Content.Find.Execute("25×78=1950", $true, $false, $false, $false, $false, $true, 1, $false, "93×81=7533", 2) | Out-Null
$d.Content.Find.Execute("31×45=1395", $true, $false, $false, $false, $false, $true, 1, $false, "40×42=1680", 2) | Out-Null
$d.Content.Find.Execute("23×32=736", $true, $false, $false, $false, $false, $true, 1, $false, "92×48=4416", 2) | Out-Null
$d.Content.Find.Execute("24×71=1704", $true, $false, $false, $false, $false, $true, 1, $false, "97×11=1067", 2) | Out-Null
$d.Content.Find.Execute("97×17=1649", $true, $false, $false, $false, $false, $true, 1, $false, "49×73=3577", 2) | Out-Null
$d.Content.Find.Execute("33×36=1188", $true, $false, $false, $false, $false, $true, 1, $false, "36×74=2664", 2) | Out-Null
$d.Content.Find.Execute("32×71=2272", $true, $false, $false, $false, $false, $true, 1, $false, "76×67=5092", 2) | Out-Null
$d.Content.Find.Execute("21×40=840", $true, $false, $false, $false, $false, $true, 1, $false, "90×85=7650", 2) | Out-Null
$d.Content.Find.Execute("51×86=4386", $true, $false, $false, $false, $false, $true, 1, $false, "28×64=1792", 2) | Out-Null
$d.Content.Find.Execute("88×53=4664", $true, $false, $false, $false, $false, $true, 1, $false, "69×94=6486", 2) | Out-Null
$d.Content.Find.Execute("34×20=680", $true, $false, $false, $false, $false, $true, 1, $false, "61×47=2867", 2) | Out-Null
$d.Content.Find.Execute("15×67=1005", $true, $false, $false, $false, $false, $true, 1, $false, "26×89=2314", 2) | Out-Null
$d.Content.Find.Execute("67×49=3283", $true, $false, $false, $false, $false, $true, 1, $false, "33×31=1023", 2) | Out-Null
$d.Content.Find.Execute("47×16=752", $true, $false, $false, $false, $false, $true, 1, $false, "27×97=2619", 2) | Out-Null
$d.Content.Find.Execute("79×98=7742", $true, $false, $false, $false, $false, $true, 1, $false, "95×46=4370", 2) | Out-Null
$d.Content.Find.Execute("82×19=1558", $true, $false, $false, $false, $false, $true, 1, $false, "54×14=756", 2) | Out-Null
$d.Content.Find.Execute("39×62=2418", $true, $false, $false, $false, $false, $true, 1, $false, "56×20=1120", 2) | Out-Null
$d.Content.Find.Execute("17×31=527", $true, $false, $false, $false, $false, $true, 1, $false, "48×84=4032", 2) | Out-Null
$d.Content.Find.Execute("64×43=2752", $true, $false, $false, $false, $false, $true, 1, $false, "15×88=1320", 2) | Out-Null
$d.Content.Find.Execute("69×73=5037", $true, $false, $false, $false, $false, $true, 1, $false, "69×37=2553", 2) | Out-Null
$d.Content.Find.Execute("35×14=490", $true, $false, $false, $false, $false, $true, 1, $false, "66×87=5742", 2) | Out-Null
$d.Content.Find.Execute("20×53=1060", $true, $false, $false, $false, $false, $true, 1, $false, "23×38=874", 2) | Out-Null
$d.Content.Find.Execute("39×34=1326", $true, $false, $false, $false, $false, $true, 1, $false, "38×45=1710", 2) | Out-Null
$d.Content.Find.Execute("20×17=340", $true, $false, $false, $false, $false, $true, 1, $false, "28×99=2772", 2) | Out-Null
